# Update the footer/version shape on the single slide of the cheatsheet deck:
#   "Package version  0.5.5 •  Updated: 2021-07"
#     -> "Package version  0.5.7 •  Updated: 2021-08"
# This corresponds to the jfa package being rebuilt at version 0.5.7.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Locate the footer shape (id 322) that holds the "Package version ... Updated: ..." text.
$shp = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $candidate = $s.Shapes.Item($i)
    if ($candidate.Id -eq 322) {
        $shp = $candidate
        break
    }
}

$tr = $shp.TextFrame.TextRange
$full = $tr.Text

# --- Bump the package version: "0.5.5" -> "0.5.7" ---------------------------
$verNeedle = "0.5.5"
$verIdx = $full.IndexOf($verNeedle)
if ($verIdx -ge 0) {
    # 1-based character position of the last "5" in "0.5.5"
    $start = $verIdx + $verNeedle.Length
    $tr.Characters($start, 1).Text = "7"
}

# Refresh full text after the first edit before locating the second needle.
$full = $tr.Text

# --- Bump the "Updated" month: "2021-07" -> "2021-08" ------------------------
$dateNeedle = "2021-07"
$dateIdx = $full.IndexOf($dateNeedle)
if ($dateIdx -ge 0) {
    # 1-based character position of the first digit of the trailing "07"
    $start = $dateIdx + ($dateNeedle.Length - 2) + 1
    $tr.Characters($start, 2).Text = "08"
}

Write-Output $tr.Text
